# "update alda 20 mei"
#
# Adds a new worksheet "pretest posttest - next page" between the existing
# "pretest posttest - search" and "manage faq - search" sheets, fills in its
# data, and updates the selected cell on all three sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Insert the new sheet right after "pretest posttest - search" ---------
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "pretest posttest - next page"

# --- Populate the new sheet's data -----------------------------------------
$newSheet.Range("A1").Value = "var_next_page"
$newSheet.Range("B1").Value = "next_page_status"
$newSheet.Range("C1").Value = "expected_output"

$newSheet.Range("B2").Value = "click"
$newSheet.Range("C2").Value = "pass"

$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = "input"
$newSheet.Range("C3").Value = "pass"

$newSheet.Range("A4").Value = 5
$newSheet.Range("B4").Value = "input"
$newSheet.Range("C4").Value = "pass"

$newSheet.Range("C5").Value = "fail"

# --- Column widths on the new sheet -----------------------------------------
$newSheet.Columns.Item(1).ColumnWidth = 17.29
$newSheet.Columns.Item(2).ColumnWidth = 16.88
$newSheet.Columns.Item(3).ColumnWidth = 17.43

# --- Update the selected cell on each sheet ---------------------------------
# Look the "manage faq" sheet up by name: its positional index shifted once
# the new sheet was inserted in front of it.
$wsFaq = $wb.Worksheets.Item("manage faq - search")

$ws1.Range("B8").Select() | Out-Null
$wsFaq.Range("B21").Select() | Out-Null

# Select last so the newly-added sheet ends up as the active tab, matching
# the workbook's activeTab.
$newSheet.Range("C10").Select() | Out-Null
